$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.080.65'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '1.660.59'
$ws.Range("E3").Value = '  -1.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.74'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5163'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.66%  '
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2579'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06290'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.95'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07527'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("D12").Value = '1.672.33'
$ws.Range("E12").Value = '  -0.79%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.402'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5374'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -5.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '66.12'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.52%  '
$ws.Range("D16").Value = '0.0₅7916'
$ws.Range("E16").Value = '  -2.78%  '
$ws.Range("D17").Value = '26.082.28'
$ws.Range("E17").Value = '  -0.62%  '
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("E19").Value = '  -3.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '187.18'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.91%  '
$ws.Range("E21").Value = '  -3.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.170'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.004'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '148.34'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1211'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -4.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.376'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.57'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.379'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06157'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -4.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.258'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.464'
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.397'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.62%  '
$ws.Range("E33").Value = '  -1.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9855'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.387'
$ws.Range("D35").ClearFormats()
$ws.Range("E36").Value = '  +1.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5865'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -4.28%  '
$ws.Range("D38").Value = '1.104.55'
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("E39").Value = '  -1.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.977'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8465'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.002'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.86'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("D44").Value = '1.809.98'
$ws.Range("E44").Value = '  -1.19%  '
$ws.Range("E45").Value = '  -0.55%  '
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.85'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.020'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05234'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4243'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.846'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.99%  '
